$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 now shows the most-recent-but-one login: username + its timestamp
$ws.Range("A1").Value = "kvw5270"
$ws.Range("B1").Value = "03/24/2020 01:10:03"

# Row 2 keeps the username and advances to the latest timestamp
$ws.Range("A2").Value = "kvw5270"
$ws.Range("B2").Value = "03/24/2020 01:12:32"

# Column A narrows slightly now that it only holds the username, not the
# wider "Username" header text.
$ws.Columns.Item(1).ColumnWidth = 8.14
